$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 114, shifting the existing rows 114-121 down to 115-122.
$ws.Rows.Item(114).Insert()

# Copy the formatting (style) of the date cell from the row above into the new row's D cell,
# so the new row matches the existing date-formatted column.
$ws.Cells.Item(113, 4).Copy()
$ws.Cells.Item(114, 4).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Populate the new row 114 with the new record's values.
$ws.Cells.Item(114, 1).Value = 11
$ws.Cells.Item(114, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(114, 3).Value = "Bíobío"
$ws.Cells.Item(114, 4).Value = 45106
$ws.Cells.Item(114, 5).Value = 8
$ws.Cells.Item(114, 6).Value = 100112037
$ws.Cells.Item(114, 7).Value = "Cebollín"
$ws.Cells.Item(114, 8).Value = "Sin especificar"
$ws.Cells.Item(114, 9).Value = "Primera"
$ws.Cells.Item(114, 10).Value = 60
$ws.Cells.Item(114, 11).Value = 7500
$ws.Cells.Item(114, 12).Value = 8000
$ws.Cells.Item(114, 13).Value = 7750
$ws.Cells.Item(114, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(114, 15).Value = "Región Metropolitana"
$ws.Cells.Item(114, 16).Value = 215
$ws.Cells.Item(114, 17).Value = 36
$ws.Cells.Item(114, 18).Value = "Hortaliza"
